# "update scripts wuth new tpm" - refresh NATMI LR-pair edge table (Nppa -> Npr1)
# with newly computed TPM-derived values.
#
# Rows 2-4 (sending cluster "ECs" -> target clusters ECs/FAPs/MuSCs) are new
# numbers for a sending cluster that was not broken out before.
# Rows 5-7 add the recomputed "MuSCs" sending-cluster rows (previously rows 2-4)
# as additional rows, now with updated TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppa"
$ws.Range("C2").Value = "Npr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1893306666666667
$ws.Range("H2").Value = 0.5679920000000001
$ws.Range("I2").Value = 0.6308162521878971
$ws.Range("J2").Value = 0.6308162521878971
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.40701066666667
$ws.Range("N2").Value = 88.22103200000001
$ws.Range("O2").Value = 0.8048351800855125
$ws.Range("P2").Value = 0.8048351800855124
$ws.Range("Q2").Value = 5.567648934193778
$ws.Range("R2").Value = 50.10884040774401
$ws.Range("S2").Value = 0.5077031119305142
$ws.Range("T2").Value = 0.5077031119305141

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppa"
$ws.Range("C3").Value = "Npr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1893306666666667
$ws.Range("H3").Value = 0.5679920000000001
$ws.Range("I3").Value = 0.6308162521878971
$ws.Range("J3").Value = 0.6308162521878971
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.578036333333333
$ws.Range("N3").Value = 19.734109
$ws.Range("O3").Value = 0.1800330920051143
$ws.Range("P3").Value = 0.1800330920051142
$ws.Range("Q3").Value = 1.245424004347556
$ws.Range("R3").Value = 11.208816039128
$ws.Range("S3").Value = 0.1135678003684651
$ws.Range("T3").Value = 0.113567800368465

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nppa"
$ws.Range("C4").Value = "Npr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1893306666666667
$ws.Range("H4").Value = 0.5679920000000001
$ws.Range("I4").Value = 0.6308162521878971
$ws.Range("J4").Value = 0.6308162521878971
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.552882
$ws.Range("N4").Value = 1.658646
$ws.Range("O4").Value = 0.0151317279093733
$ws.Range("P4").Value = 0.01513172790937329
$ws.Range("Q4").Value = 0.104677517648
$ws.Range("R4").Value = 0.942097658832
$ws.Range("S4").Value = 0.009545339888917866
$ws.Range("T4").Value = 0.009545339888917865

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nppa"
$ws.Range("C5").Value = "Npr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1108053333333333
$ws.Range("H5").Value = 0.332416
$ws.Range("I5").Value = 0.3691837478121029
$ws.Range("J5").Value = 0.3691837478121029
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 29.40701066666667
$ws.Range("N5").Value = 88.22103200000001
$ws.Range("O5").Value = 0.8048351800855125
$ws.Range("P5").Value = 0.8048351800855124
$ws.Range("Q5").Value = 3.258453619256889
$ws.Range("R5").Value = 29.326082573312
$ws.Range("S5").Value = 0.2971320681549983
$ws.Range("T5").Value = 0.2971320681549982

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nppa"
$ws.Range("C6").Value = "Npr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1108053333333333
$ws.Range("H6").Value = 0.332416
$ws.Range("I6").Value = 0.3691837478121029
$ws.Range("J6").Value = 0.3691837478121029
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.578036333333333
$ws.Range("N6").Value = 19.734109
$ws.Range("O6").Value = 0.1800330920051143
$ws.Range("P6").Value = 0.1800330920051142
$ws.Range("Q6").Value = 0.7288815085937778
$ws.Range("R6").Value = 6.559933577344
$ws.Range("S6").Value = 0.06646529163664924
$ws.Range("T6").Value = 0.06646529163664922

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Nppa"
$ws.Range("C7").Value = "Npr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1108053333333333
$ws.Range("H7").Value = 0.332416
$ws.Range("I7").Value = 0.3691837478121029
$ws.Range("J7").Value = 0.3691837478121029
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.552882
$ws.Range("N7").Value = 1.658646
$ws.Range("O7").Value = 0.0151317279093733
$ws.Range("P7").Value = 0.01513172790937329
$ws.Range("Q7").Value = 0.06126227430399999
$ws.Range("R7").Value = 0.551360468736
$ws.Range("S7").Value = 0.005586388020455431
$ws.Range("T7").Value = 0.005586388020455429
